$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1357.0303
$ws.Range("I19").Value = 1798.9375
$ws.Range("J19").Value = 941.1177
$ws.Range("K19").Value = 1798.9375
$ws.Range("L19").Value = 941.1177
$ws.Range("M19").Value = -1623.9375
$ws.Range("N19").Value = -1291.1177

$ws.Range("H70").Value = 2891.6667
$ws.Range("J70").Value = 2891.6667
$ws.Range("L70").Value = 8675.000100000001
$ws.Range("N70").Value = -9215.000100000001

$ws.Range("H73").Value = 2891.6667
$ws.Range("J73").Value = 2891.6667
$ws.Range("L73").Value = 8675.000100000001
$ws.Range("N73").Value = -10547.0001

$ws.Range("H107").Value = 1715.25
$ws.Range("I107").Value = 1911.8572
$ws.Range("J107").Value = 1440
$ws.Range("K107").Value = 1911.8572
$ws.Range("L107").Value = 1440
$ws.Range("M107").Value = 8.142800000000079
$ws.Range("N107").Value = -5280

$ws.Range("H113").Value = 1301.6
$ws.Range("I113").Value = 755
$ws.Range("J113").Value = 1438.25
$ws.Range("K113").Value = 755
$ws.Range("L113").Value = 1438.25
$ws.Range("M113").Value = 2499
$ws.Range("N113").Value = -7946.25

$ws.Range("H132").Value = 155677.05
$ws.Range("I132").Value = 163158.19
$ws.Range("J132").Value = 1066.6666
$ws.Range("K132").Value = 489474.57
$ws.Range("L132").Value = 3199.9998
$ws.Range("M132").Value = -486944.57
$ws.Range("N132").Value = -8259.9998

$ws.Range("H138").Value = 1889084
$ws.Range("I138").Value = 2779459.2
$ws.Range("J138").Value = 3583.4119
$ws.Range("K138").Value = 8338377.600000001
$ws.Range("L138").Value = 10750.2357
$ws.Range("M138").Value = -8333237.600000001
$ws.Range("N138").Value = -21030.2357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5283.3965
$ws.Range("I32").Value = 5578.8086
$ws.Range("J32").Value = 4021.182
$ws.Range("K32").Value = 5578.8086
$ws.Range("L32").Value = 4021.182
$ws.Range("M32").Value = -5291.8086
$ws.Range("N32").Value = -4595.182

$ws.Range("H45").Value = 1075.0769
$ws.Range("I45").Value = 970.8570999999999
$ws.Range("J45").Value = 1196.6666
$ws.Range("K45").Value = 970.8570999999999
$ws.Range("L45").Value = 1196.6666
$ws.Range("M45").Value = -593.8570999999999
$ws.Range("N45").Value = -1950.6666

$ws.Range("H74").Value = 4028.318
$ws.Range("I74").Value = 876.03125
$ws.Range("J74").Value = 12434.417
$ws.Range("K74").Value = 876.03125
$ws.Range("L74").Value = 12434.417
$ws.Range("M74").Value = -2.03125
$ws.Range("N74").Value = -14182.417

$ws.Range("H77").Value = 4028.318
$ws.Range("I77").Value = 876.03125
$ws.Range("J77").Value = 12434.417
$ws.Range("K77").Value = 4380.15625
$ws.Range("L77").Value = 62172.085
$ws.Range("M77").Value = -12.15625
$ws.Range("N77").Value = -70908.08499999999

$ws.Range("H122").Value = 1601
$ws.Range("I122").Value = 1671.2
$ws.Range("J122").Value = 1250
$ws.Range("K122").Value = 5013.6
$ws.Range("L122").Value = 3750
$ws.Range("M122").Value = -2563.6
$ws.Range("N122").Value = -8650

$ws.Range("H132").Value = 34614.13
$ws.Range("I132").Value = 2351.0417
$ws.Range("J132").Value = 145230.42
$ws.Range("K132").Value = 7053.125100000001
$ws.Range("L132").Value = 435691.26
$ws.Range("M132").Value = -4523.125100000001
$ws.Range("N132").Value = -440751.26

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 59592.42
$ws.Range("I134").Value = 66309.17999999999
$ws.Range("K134").Value = 198927.54
$ws.Range("M134").Value = -196392.54

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1446.909
$ws.Range("I16").Value = 1648.3334
$ws.Range("J16").Value = 1205.2
$ws.Range("K16").Value = 1648.3334
$ws.Range("L16").Value = 1205.2
$ws.Range("M16").Value = -1361.3334
$ws.Range("N16").Value = -1779.2

$ws.Range("H31").Value = 1326.814
$ws.Range("I31").Value = 1044.9166
$ws.Range("J31").Value = 1682.8948
$ws.Range("K31").Value = 1044.9166
$ws.Range("L31").Value = 1682.8948
$ws.Range("M31").Value = -749.9166
$ws.Range("N31").Value = -2272.8948

$ws.Range("H34").Value = 1326.814
$ws.Range("I34").Value = 1044.9166
$ws.Range("J34").Value = 1682.8948
$ws.Range("K34").Value = 1044.9166
$ws.Range("L34").Value = 1682.8948
$ws.Range("M34").Value = -842.9166
$ws.Range("N34").Value = -2086.8948

$ws.Range("H113").Value = 1446.909
$ws.Range("I113").Value = 1648.3334
$ws.Range("J113").Value = 1205.2
$ws.Range("K113").Value = 1648.3334
$ws.Range("L113").Value = 1205.2
$ws.Range("M113").Value = 521.6666
$ws.Range("N113").Value = -5545.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 610.5
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 610.5
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 1831.5
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -2055.5

$ws.Range("H34").Value = 62503076
$ws.Range("I34").Value = 333.75
$ws.Range("J34").Value = 83337320
$ws.Range("K34").Value = 1001.25
$ws.Range("L34").Value = 250011960
$ws.Range("M34").Value = -917.25
$ws.Range("N34").Value = -250012128

$ws.Range("H39").Value = 2971.75
$ws.Range("I39").Value = 993
$ws.Range("J39").Value = 3631.3333
$ws.Range("K39").Value = 2979
$ws.Range("L39").Value = 10893.9999
$ws.Range("M39").Value = -2685
$ws.Range("N39").Value = -11481.9999

$ws.Range("H80").Value = 7244112.5
$ws.Range("J80").Value = 2168
$ws.Range("L80").Value = 6504
$ws.Range("N80").Value = -8376

$ws.Range("H83").Value = 7244112.5
$ws.Range("J83").Value = 2168
$ws.Range("L83").Value = 19512
$ws.Range("N83").Value = -28872

$ws.Range("H122").Value = 19609840
$ws.Range("I122").Value = 41667404
$ws.Range("J122").Value = 3116.3333
$ws.Range("K122").Value = 375006636
$ws.Range("L122").Value = 28046.9997
$ws.Range("M122").Value = -375004186
$ws.Range("N122").Value = -32946.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 406.72415
$ws.Range("I107").Value = 252.82353
$ws.Range("J107").Value = 624.75
$ws.Range("K107").Value = 252.82353
$ws.Range("L107").Value = 624.75
$ws.Range("M107").Value = 1667.17647
$ws.Range("N107").Value = -4464.75

$ws.Range("H132").Value = 1937.5385
$ws.Range("I132").Value = 1852.1538
$ws.Range("J132").Value = 2108.3076
$ws.Range("K132").Value = 5556.4614
$ws.Range("L132").Value = 6324.9228
$ws.Range("M132").Value = -3026.4614
$ws.Range("N132").Value = -11384.9228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 568.63635
$ws.Range("I22").Value = 570.25
$ws.Range("J22").Value = 567.7143
$ws.Range("K22").Value = 570.25
$ws.Range("L22").Value = 567.7143
$ws.Range("M22").Value = -275.25
$ws.Range("N22").Value = -1157.7143

$ws.Range("H27").Value = 568.63635
$ws.Range("I27").Value = 570.25
$ws.Range("J27").Value = 567.7143
$ws.Range("K27").Value = 570.25
$ws.Range("L27").Value = 567.7143
$ws.Range("M27").Value = -463.25
$ws.Range("N27").Value = -781.7143

$ws.Range("H40").Value = 1420
$ws.Range("I40").Value = 1420
$ws.Range("K40").Value = 1420
$ws.Range("M40").Value = -1284

$ws.Range("H93").Value = 5500
$ws.Range("I93").Value = 6500
$ws.Range("J93").Value = 3500
$ws.Range("K93").Value = 6500
$ws.Range("L93").Value = 3500
$ws.Range("M93").Value = -5252
$ws.Range("N93").Value = -5996

$ws.Range("H122").Value = 2086.1428
$ws.Range("I122").Value = 1302
$ws.Range("J122").Value = 2399.8
$ws.Range("K122").Value = 3906
$ws.Range("L122").Value = 7199.400000000001
$ws.Range("M122").Value = -1456
$ws.Range("N122").Value = -12099.4

$ws.Range("H132").Value = 4247.75
$ws.Range("I132").Value = 5087.8335
$ws.Range("J132").Value = 2987.625
$ws.Range("K132").Value = 15263.5005
$ws.Range("L132").Value = 8962.875
$ws.Range("M132").Value = -12733.5005
$ws.Range("N132").Value = -14022.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 480.14285
$ws.Range("I113").Value = 535
$ws.Range("J113").Value = 439
$ws.Range("K113").Value = 1605
$ws.Range("L113").Value = 1317
$ws.Range("M113").Value = 565
$ws.Range("N113").Value = -5657

$ws.Range("H132").Value = 4585.9775
$ws.Range("I132").Value = 5384
$ws.Range("K132").Value = 16152
$ws.Range("M132").Value = -13622

